$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update "总计" (sheet1): insert a new row for 2022-Q4 data
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").Style = "Normal"
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 24
$total.Cells.Item(2, 4).Value = 6.32

# copy the numeric-index cell style (s="2") from the row below onto the new row
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)

# fix the running index values (column A) for the rows that got shifted down
for ($r = 3; $r -le 8; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) Insert a brand new worksheet "2022-Q4" right after "总计"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# apply the shared header style (s="2", bold + border + center/top alignment)
# used on every other sheet in this workbook
$srcHeader = $wb.Worksheets.Item("2022-Q3").Range("B1:H1")
$srcHeader.Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$data = @(
    @('515220','国泰中证煤炭ETF','45.22','99.89','4.59','2.0756','8'),
    @('012526','广发盛锦混合A','24.00','93.13','3.81','0.9144','9'),
    @('213001','宝盈鸿利收益灵活配置混合A','14.47','91.18','5.57','0.8060','3'),
    @('161032','富国中证煤炭指数A','16.36','94.14','4.34','0.7100','8'),
    @('013275','富国中证煤炭指数C','12.51','94.14','4.34','0.5429','8'),
    @('168204','中融中证煤炭指数A','8.35','91.94','4.22','0.3524','8'),
    @('011336','兴全汇吉一年持有期混合A','15.09','39.83','2.19','0.3305','4'),
    @('630011','华商主题精选混合','3.43','92.61','3.84','0.1317','9'),
    @('001543','宝盈新锐灵活配置混合A','2.28','91.59','4.86','0.1108','8'),
    @('159930','汇添富中证能源ETF','2.14','99.31','3.44','0.0736','9'),
    @('009965','宝盈祥琪混合A','1.40','43.98','3.44','0.0482','5'),
    @('012527','广发盛锦混合C','1.14','93.13','3.81','0.0434','9'),
    @('007581','宝盈鸿利收益灵活配置混合C','0.62','91.18','5.57','0.0345','3'),
    @('011997','景顺长城安盈回报一年持有期混合A','1.50','26.78','1.88','0.0282','2'),
    @('260117','景顺长城支柱产业混合','0.77','72.88','3.23','0.0249','10'),
    @('014768','景顺华城稳健6月持有混合C','1.61','22.61','1.41','0.0227','3'),
    @('011337','兴全汇吉一年持有期混合C','0.80','39.83','2.19','0.0175','4'),
    @('014767','景顺华城稳健6月持有混合A','1.10','22.61','1.41','0.0155','3'),
    @('007578','宝盈新锐灵活配置混合C','0.21','91.59','4.86','0.0102','8'),
    @('016814','中融中证煤炭指数C','0.24','91.94','4.22','0.0101','8'),
    @('001135','益民品质升级灵活配置混合','0.47','74.18','1.65','0.0078','9'),
    @('008890','中邮价值优选一年定期开放灵活配置混合','0.12','64.16','3.94','0.0047','5'),
    @('011998','景顺长城安盈回报一年持有期混合C','0.08','26.78','1.88','0.0015','2'),
    @('009966','宝盈祥琪混合C','0.01','43.98','3.44','0.0003','5')

)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]

    $q4.Cells.Item($row, 1).Value = $i

    $q4.Cells.Item($row, 2).NumberFormat = "@"
    $q4.Cells.Item($row, 2).Value = $vals[0]

    $q4.Cells.Item($row, 3).NumberFormat = "@"
    $q4.Cells.Item($row, 3).Value = $vals[1]

    $q4.Cells.Item($row, 4).NumberFormat = "@"
    $q4.Cells.Item($row, 4).Value = $vals[2]

    $q4.Cells.Item($row, 5).NumberFormat = "@"
    $q4.Cells.Item($row, 5).Value = $vals[3]

    $q4.Cells.Item($row, 6).NumberFormat = "@"
    $q4.Cells.Item($row, 6).Value = $vals[4]

    $q4.Cells.Item($row, 7).NumberFormat = "@"
    $q4.Cells.Item($row, 7).Value = $vals[5]

    $q4.Cells.Item($row, 8).Value = [int]$vals[6]
}

# apply the shared numeric-index style (s="2") from the header row onto all of
# column A's data cells, then re-write the index values (paste special resets them)
$q4.Range("B1").Copy()
$q4.Range("A2:A" + (1 + $data.Length)).PasteSpecial(-4122)
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $q4.Cells.Item($row, 1).Value = $i
}
